# "files added for AK Jain" -- populate the "Personal Web" sub-columns
# (H, K, N, Q, T, W, Z, AC, AF, AI of the per-source-group triples) for
# Anil K. Jain (row 7), add the column-E "Personal Web total" formula for
# every data row (it references those sub-columns), and restore the
# scrolled/selected view state.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 5 sub-headers: the "Personal Web" column of each of the 10 source
# groups (G:I, J:L, ..., AH:AJ) was blank; it now carries the same
# "Personal Web" label used elsewhere in the header row (shared string).
$personalWebCols = @("H","K","N","Q","T","W","Z","AC","AF","AI")
foreach ($col in $personalWebCols) {
    $ws.Range($col + "5").Value = "Personal Web"
}

# Row 7 (Anil K. Jain) - new "Personal Web" publication counts per source.
$ws.Range("H7").Value  = 154
$ws.Range("K7").Value  = 276
$ws.Range("N7").Value  = 145
$ws.Range("Q7").Value  = 57
$ws.Range("T7").Value  = 34
$ws.Range("W7").Value  = 5
$ws.Range("Z7").Value  = 1
$ws.Range("AC7").Value = 0
$ws.Range("AF7").Value = 0
$ws.Range("AI7").Value = 1

# Column E ("Personal Web" total) formula, added for every author row
# (6 already had it; 7-55 gain it now). For rows other than 7 the new
# H/K/N/Q/T/W/Z/AC/AF/AI inputs are still blank, so the formula evaluates
# to 0 there.
for ($row = 7; $row -le 55; $row++) {
    $ws.Range("E$row").Formula = "=H$row+K$row+N$row+Q$row+T$row+W$row+Z$row+AC$row+AF$row+AI$row"
}

# Restore the saved view state: scrolled so row 58 is the top visible
# row, with A6 selected/active.
$ws.Range("A6").Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 58
$win.ScrollColumn = 1
